$d = $word.ActiveDocument

# Locate the paragraph that holds the "Que:03" heading text.
$rng = $d.Content
$rng.Find.Execute("Que:03", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)
$r = $para.Range

# Build a WordOpenXML fragment for just this paragraph that preserves all of
# its existing content/attributes but adds bold (w:b / w:bCs) formatting to
# the paragraph mark and to every run in the paragraph, matching what Word
# itself produces when the whole paragraph (including its pilcrow) is
# selected and Bold is toggled on.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14">' +
          '<w:body>' +
            '<w:p w14:paraId="62F77BA4" w14:textId="3068C2E0" w:rsidR="00F731DB" w:rsidRDefault="00702520">' +
              '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Que:0</w:t></w:r>' +
              '<w:r w:rsidR="00867A5E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r.InsertXML($xml) | Out-Null
